$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "03/15-03/28" paragraph: append a trailing space run
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(4)
$r = $p.Range
$insertPos = $r.End - 1
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 2) "How many hours you have to work on this project this sprint" paragraph:
#    remove the grammar-check proofErr markers around "have to" by deleting
#    the text (not the paragraph mark) and retyping as one run, then append a
#    trailing space run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(5)
$r = $p.Range
$full = $d.Range($r.Start, $r.End - 1)
$full.Delete()
$ins = $d.Paragraphs.Item(5).Range
$ins.Collapse(1)
$ins.InsertAfter("How many hours you have to work on this project this sprint")
$r5 = $d.Paragraphs.Item(5).Range
$insertPos = $r5.End - 1
$sp = $d.Range($insertPos, $insertPos)
$sp.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 3) "Alex - somewhat busy, can allocate 8-12hrs over sprint" ->
#    "Alex - somewhat busy, can allocate 10-14hrs over sprint"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("8-12hrs", $true, $false, $false, $false, $false, $true, 1, $false, "10-14hrs", 2)

# ---------------------------------------------------------------------------
# 4) "Who is going to be here on what days (vacation / other class
#    priorities)" paragraph: remove the grammar-check proofErr markers. The
#    trailing proofErr sits right against the paragraph mark, so the text-only
#    delete used above would leave it behind; instead delete the whole
#    paragraph range (mark included), retype it, and then explicitly restore
#    the paragraph formatting on both halves of the old mark-merge so the
#    list level stays correct. Finally append a trailing space run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$full = $d.Range($r.Start, $r.End)
$full.Delete()
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertAfter("Who is going to be here on what days (vacation / other class priorities)`r")

$p9 = $d.Paragraphs.Item(9)
$p9.Style = "List Paragraph"
$p9.Range.ListFormat.ListLevelNumber = 1

$p10 = $d.Paragraphs.Item(10)
$p10.Style = "List Paragraph"
$p10.Range.ListFormat.ListLevelNumber = 2

$r9 = $d.Paragraphs.Item(9).Range
$insertPos = $r9.End - 1
$sp = $d.Range($insertPos, $insertPos)
$sp.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 5) "Pick a metric to evaluate each other on during the retrospective"
#    paragraph: append a space, then a gray "-same metric for every sprint"
#    run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$r = $p.Range
$insertPos = $r.End - 1
$sp = $d.Range($insertPos, $insertPos)
$sp.InsertAfter(" ")

$r13 = $d.Paragraphs.Item(13).Range
$insertPos2 = $r13.End - 1
$tail = $d.Range($insertPos2, $insertPos2)
$tail.InsertAfter("-same metric for every sprint")
$tailLen = "-same metric for every sprint".Length
$colorRange = $d.Range($insertPos2, $insertPos2 + $tailLen)
$colorRange.Font.Color = 7434614

# ---------------------------------------------------------------------------
# 6) "Communication concerning implementation needs and needed areas of
#    focus" paragraph: append ". Individual difficulties"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(16)
$r = $p.Range
$insertPos = $r.End - 1
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter(". Individual difficulties")

# ---------------------------------------------------------------------------
# 7) Remove the "Fill out details for each story" ... "Attach a screenshot of
#    your Sprint Backlog after planning" section entirely (paragraphs 19-37),
#    leaving "Alex Thurston - Scrum Master for Sprint 2" followed directly by
#    the trailing empty paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(19)
$pEnd = $d.Paragraphs.Item(37)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
